$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: usuario_asignado ---
$ws.Range("G1").Value = "usuario_asignado"

# Apply Text number format to the existing D-column id cells (D2:D4) without
# disturbing their numeric storage.
$ws.Range("D2:D4").NumberFormat = "@"

# Numeric "weights" for the first three existing rows - set the values BEFORE
# applying the Text format so they stay stored as numbers (matches target
# <v>1.2</v> / <v>3.4</v> rather than being coerced into shared strings).
$ws.Range("G2").Value = 1.2
$ws.Range("G3").Value = 3.4
$ws.Range("G2:G4").NumberFormat = "@"
# G4 intentionally left blank (formatted only).

# --- New row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 465
$ws.Range("C5").Value = 22
$ws.Range("E5").Value = "Otro centro"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245"
$ws.Range("F5").Value = "Calle X con carrera Y"

# --- New row 6 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 465
$ws.Range("C6").Value = 22
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "258"
$ws.Range("E6").Value = "Otro lado"
$ws.Range("F6").Value = "Calle X con carrera Y"

# New column G gets its own best-fit-ish width, matching the other
# "*_id"-sized text columns (e.g. column C / departamento_id).
$ws.Columns("G").ColumnWidth = $ws.Columns("C").ColumnWidth

# Match the author's final selection / active cell.
$ws.Range("G5").Select() | Out-Null
